$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 76
$ws.Cells.Item(76, 8).Value2 = 1000  # ALC!H76
$ws.Cells.Item(76, 9).Value2 = 0  # ALC!I76
$ws.Cells.Item(76, 11).Value2 = 0  # ALC!K76
$ws.Cells.Item(76, 13).ClearContents()  # ALC!M76

# Row 79
$ws.Cells.Item(79, 8).Value2 = 1000  # ALC!H79
$ws.Cells.Item(79, 9).Value2 = 0  # ALC!I79
$ws.Cells.Item(79, 11).Value2 = 0  # ALC!K79
$ws.Cells.Item(79, 13).ClearContents()  # ALC!M79

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Cells.Item(2, 8).Value2 = 1936.3636  # ARM!H2
$ws.Cells.Item(2, 9).Value2 = 825.25  # ARM!I2
$ws.Cells.Item(2, 10).Value2 = 4899.3335  # ARM!J2
$ws.Cells.Item(2, 11).Value2 = 825.25  # ARM!K2
$ws.Cells.Item(2, 12).Value2 = 4899.3335  # ARM!L2
$ws.Cells.Item(2, 13).Value2 = -712.25  # ARM!M2
$ws.Cells.Item(2, 14).Value2 = -5125.3335  # ARM!N2

# Row 32
$ws.Cells.Item(32, 8).Value2 = 6900.619  # ARM!H32
$ws.Cells.Item(32, 9).Value2 = 4723.0835  # ARM!I32
$ws.Cells.Item(32, 11).Value2 = 4723.0835  # ARM!K32
$ws.Cells.Item(32, 13).Value2 = -4436.0835  # ARM!M32

# Row 45
$ws.Cells.Item(45, 8).Value2 = 4498  # ARM!H45
$ws.Cells.Item(45, 9).Value2 = 4498  # ARM!I45
$ws.Cells.Item(45, 11).Value2 = 4498  # ARM!K45
$ws.Cells.Item(45, 13).Value2 = -4121  # ARM!M45

# Row 102
$ws.Cells.Item(102, 8).Value2 = 1833.8  # ARM!H102
$ws.Cells.Item(102, 9).Value2 = 390  # ARM!I102
$ws.Cells.Item(102, 10).Value2 = 3999.5  # ARM!J102
$ws.Cells.Item(102, 11).Value2 = 390  # ARM!K102
$ws.Cells.Item(102, 12).Value2 = 3999.5  # ARM!L102
$ws.Cells.Item(102, 13).Value2 = 1232  # ARM!M102
$ws.Cells.Item(102, 14).Value2 = -7243.5  # ARM!N102

# Row 116
$ws.Cells.Item(116, 8).Value2 = 1936.3636  # ARM!H116
$ws.Cells.Item(116, 9).Value2 = 825.25  # ARM!I116
$ws.Cells.Item(116, 10).Value2 = 4899.3335  # ARM!J116
$ws.Cells.Item(116, 11).Value2 = 825.25  # ARM!K116
$ws.Cells.Item(116, 12).Value2 = 4899.3335  # ARM!L116
$ws.Cells.Item(116, 13).Value2 = 1468.75  # ARM!M116
$ws.Cells.Item(116, 14).Value2 = -9487.333500000001  # ARM!N116

# Row 122
$ws.Cells.Item(122, 8).Value2 = 2991.2144  # ARM!H122
$ws.Cells.Item(122, 9).Value2 = 2573.1667  # ARM!I122
$ws.Cells.Item(122, 11).Value2 = 7719.500100000001  # ARM!K122
$ws.Cells.Item(122, 13).Value2 = -5269.500100000001  # ARM!M122

# Row 138
$ws.Cells.Item(138, 8).Value2 = 74756.336  # ARM!H138
$ws.Cells.Item(138, 10).Value2 = 74756.336  # ARM!J138
$ws.Cells.Item(138, 12).Value2 = 74756.336  # ARM!L138
$ws.Cells.Item(138, 14).Value2 = -85036.336  # ARM!N138

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Cells.Item(3, 8).Value2 = 1936.3636  # BSM!H3
$ws.Cells.Item(3, 9).Value2 = 825.25  # BSM!I3
$ws.Cells.Item(3, 10).Value2 = 4899.3335  # BSM!J3
$ws.Cells.Item(3, 11).Value2 = 825.25  # BSM!K3
$ws.Cells.Item(3, 12).Value2 = 4899.3335  # BSM!L3
$ws.Cells.Item(3, 13).Value2 = -711.25  # BSM!M3
$ws.Cells.Item(3, 14).Value2 = -5127.3335  # BSM!N3

# Row 105
$ws.Cells.Item(105, 8).Value2 = 2998.3572  # BSM!H105
$ws.Cells.Item(105, 9).Value2 = 2748.1667  # BSM!I105
$ws.Cells.Item(105, 11).Value2 = 2748.1667  # BSM!K105
$ws.Cells.Item(105, 13).Value2 = -1001.1667  # BSM!M105

# Row 107
$ws.Cells.Item(107, 8).Value2 = 8400  # BSM!H107
$ws.Cells.Item(107, 9).Value2 = 8500  # BSM!I107
$ws.Cells.Item(107, 10).Value2 = 8000  # BSM!J107
$ws.Cells.Item(107, 11).Value2 = 8500  # BSM!K107
$ws.Cells.Item(107, 12).Value2 = 8000  # BSM!L107
$ws.Cells.Item(107, 13).Value2 = -6580  # BSM!M107
$ws.Cells.Item(107, 14).Value2 = -11840  # BSM!N107

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Cells.Item(22, 8).Value2 = 384.13333  # CRP!H22
$ws.Cells.Item(22, 10).Value2 = 405.15384  # CRP!J22
$ws.Cells.Item(22, 12).Value2 = 405.15384  # CRP!L22
$ws.Cells.Item(22, 14).Value2 = -1105.15384  # CRP!N22

# Row 97
$ws.Cells.Item(97, 8).Value2 = 16000  # CRP!H97
$ws.Cells.Item(97, 9).Value2 = 8000  # CRP!I97
$ws.Cells.Item(97, 10).Value2 = 24000  # CRP!J97
$ws.Cells.Item(97, 11).Value2 = 8000  # CRP!K97
$ws.Cells.Item(97, 12).Value2 = 24000  # CRP!L97
$ws.Cells.Item(97, 13).Value2 = -7009  # CRP!M97
$ws.Cells.Item(97, 14).Value2 = -25982  # CRP!N97

# Row 99
$ws.Cells.Item(99, 8).Value2 = 14476.143  # CRP!H99
$ws.Cells.Item(99, 9).Value2 = 12674.143  # CRP!I99
$ws.Cells.Item(99, 11).Value2 = 12674.143  # CRP!K99
$ws.Cells.Item(99, 13).Value2 = -11176.143  # CRP!M99

# Row 122
$ws.Cells.Item(122, 8).Value2 = 2644.0527  # CRP!H122
$ws.Cells.Item(122, 9).Value2 = 2652.0557  # CRP!I122
$ws.Cells.Item(122, 10).Value2 = 2500  # CRP!J122
$ws.Cells.Item(122, 11).Value2 = 7956.1671  # CRP!K122
$ws.Cells.Item(122, 12).Value2 = 7500  # CRP!L122
$ws.Cells.Item(122, 13).Value2 = -5506.1671  # CRP!M122
$ws.Cells.Item(122, 14).Value2 = -12400  # CRP!N122

# Row 126
$ws.Cells.Item(126, 8).Value2 = 14476.143  # CRP!H126
$ws.Cells.Item(126, 9).Value2 = 12674.143  # CRP!I126
$ws.Cells.Item(126, 11).Value2 = 38022.429  # CRP!K126
$ws.Cells.Item(126, 13).Value2 = -35552.429  # CRP!M126

# Row 130
$ws.Cells.Item(130, 8).Value2 = 54250  # CRP!H130
$ws.Cells.Item(130, 9).Value2 = 49000  # CRP!I130
$ws.Cells.Item(130, 11).Value2 = 49000  # CRP!K130
$ws.Cells.Item(130, 13).Value2 = -43980  # CRP!M130

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 38
$ws.Cells.Item(38, 8).Value2 = 82.083336  # CUL!H38
$ws.Cells.Item(38, 9).Value2 = 29.428572  # CUL!I38
$ws.Cells.Item(38, 10).Value2 = 155.8  # CUL!J38
$ws.Cells.Item(38, 11).Value2 = 88.28571599999999  # CUL!K38
$ws.Cells.Item(38, 12).Value2 = 467.4  # CUL!L38
$ws.Cells.Item(38, 13).Value2 = 258.714284  # CUL!M38
$ws.Cells.Item(38, 14).Value2 = -1161.4  # CUL!N38

# Row 131
$ws.Cells.Item(131, 8).Value2 = 1331.3334  # CUL!H131
$ws.Cells.Item(131, 10).Value2 = 2665  # CUL!J131
$ws.Cells.Item(131, 12).Value2 = 7995  # CUL!L131
$ws.Cells.Item(131, 14).Value2 = -18075  # CUL!N131

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Cells.Item(102, 8).Value2 = 2138.111  # GSM!H102
$ws.Cells.Item(102, 9).Value2 = 2138.111  # GSM!I102
$ws.Cells.Item(102, 11).Value2 = 2138.111  # GSM!K102
$ws.Cells.Item(102, 13).Value2 = -516.1109999999999  # GSM!M102

# Row 113
$ws.Cells.Item(113, 9).Value2 = 1149.5  # GSM!I113
$ws.Cells.Item(113, 10).Value2 = 4250  # GSM!J113
$ws.Cells.Item(113, 11).Value2 = 1149.5  # GSM!K113
$ws.Cells.Item(113, 12).Value2 = 4250  # GSM!L113
$ws.Cells.Item(113, 13).Value2 = 1020.5  # GSM!M113
$ws.Cells.Item(113, 14).Value2 = -8590  # GSM!N113

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Cells.Item(7, 8).Value2 = 4538.3335  # LTW!H7
$ws.Cells.Item(7, 9).Value2 = 4538.3335  # LTW!I7
$ws.Cells.Item(7, 11).Value2 = 4538.3335  # LTW!K7
$ws.Cells.Item(7, 13).Value2 = -4426.3335  # LTW!M7

# Row 22
$ws.Cells.Item(22, 8).Value2 = 3489.2068  # LTW!H22
$ws.Cells.Item(22, 9).Value2 = 2453.8262  # LTW!I22
$ws.Cells.Item(22, 10).Value2 = 7458.1665  # LTW!J22
$ws.Cells.Item(22, 11).Value2 = 2453.8262  # LTW!K22
$ws.Cells.Item(22, 12).Value2 = 7458.1665  # LTW!L22
$ws.Cells.Item(22, 13).Value2 = -2158.8262  # LTW!M22
$ws.Cells.Item(22, 14).Value2 = -8048.1665  # LTW!N22

# Row 27
$ws.Cells.Item(27, 8).Value2 = 3489.2068  # LTW!H27
$ws.Cells.Item(27, 9).Value2 = 2453.8262  # LTW!I27
$ws.Cells.Item(27, 10).Value2 = 7458.1665  # LTW!J27
$ws.Cells.Item(27, 11).Value2 = 2453.8262  # LTW!K27
$ws.Cells.Item(27, 12).Value2 = 7458.1665  # LTW!L27
$ws.Cells.Item(27, 13).Value2 = -2346.8262  # LTW!M27
$ws.Cells.Item(27, 14).Value2 = -7672.1665  # LTW!N27

# Row 40
$ws.Cells.Item(40, 8).Value2 = 2598.8462  # LTW!H40
$ws.Cells.Item(40, 9).Value2 = 2709  # LTW!I40
$ws.Cells.Item(40, 10).Value2 = 1993  # LTW!J40
$ws.Cells.Item(40, 11).Value2 = 2709  # LTW!K40
$ws.Cells.Item(40, 12).Value2 = 1993  # LTW!L40
$ws.Cells.Item(40, 13).Value2 = -2573  # LTW!M40
$ws.Cells.Item(40, 14).Value2 = -2265  # LTW!N40

# Row 46
$ws.Cells.Item(46, 8).Value2 = 4767.4  # LTW!H46
$ws.Cells.Item(46, 9).Value2 = 2499.5  # LTW!I46
$ws.Cells.Item(46, 10).Value2 = 6279.3335  # LTW!J46
$ws.Cells.Item(46, 11).Value2 = 2499.5  # LTW!K46
$ws.Cells.Item(46, 12).Value2 = 6279.3335  # LTW!L46
$ws.Cells.Item(46, 13).Value2 = -2311.5  # LTW!M46
$ws.Cells.Item(46, 14).Value2 = -6655.3335  # LTW!N46

# Row 55
$ws.Cells.Item(55, 8).Value2 = 811.7692  # LTW!H55
$ws.Cells.Item(55, 9).Value2 = 570.5  # LTW!I55
$ws.Cells.Item(55, 11).Value2 = 570.5  # LTW!K55
$ws.Cells.Item(55, 13).Value2 = -397.5  # LTW!M55

# Row 61
$ws.Cells.Item(61, 8).Value2 = 3640.2222  # LTW!H61
$ws.Cells.Item(61, 9).Value2 = 3640.2222  # LTW!I61
$ws.Cells.Item(61, 10).Value2 = 0  # LTW!J61
$ws.Cells.Item(61, 11).Value2 = 3640.2222  # LTW!K61
$ws.Cells.Item(61, 12).Value2 = 0  # LTW!L61
$ws.Cells.Item(61, 13).Value2 = -3438.2222  # LTW!M61
$ws.Cells.Item(61, 14).ClearContents()  # LTW!N61

# Row 113
$ws.Cells.Item(113, 8).Value2 = 3640.2222  # LTW!H113
$ws.Cells.Item(113, 9).Value2 = 3640.2222  # LTW!I113
$ws.Cells.Item(113, 10).Value2 = 0  # LTW!J113
$ws.Cells.Item(113, 11).Value2 = 3640.2222  # LTW!K113
$ws.Cells.Item(113, 12).Value2 = 0  # LTW!L113
$ws.Cells.Item(113, 13).Value2 = -1470.2222  # LTW!M113
$ws.Cells.Item(113, 14).ClearContents()  # LTW!N113

# Row 122
$ws.Cells.Item(122, 8).Value2 = 3397.8  # LTW!H122
$ws.Cells.Item(122, 9).Value2 = 3397.8  # LTW!I122
$ws.Cells.Item(122, 11).Value2 = 10193.4  # LTW!K122
$ws.Cells.Item(122, 13).Value2 = -7743.400000000001  # LTW!M122

# Row 126
$ws.Cells.Item(126, 8).Value2 = 4538.3335  # LTW!H126
$ws.Cells.Item(126, 9).Value2 = 4538.3335  # LTW!I126
$ws.Cells.Item(126, 11).Value2 = 13615.0005  # LTW!K126
$ws.Cells.Item(126, 13).Value2 = -11145.0005  # LTW!M126

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Cells.Item(107, 8).Value2 = 1057.5  # WVR!H107
$ws.Cells.Item(107, 9).Value2 = 1062.7273  # WVR!I107
$ws.Cells.Item(107, 11).Value2 = 3188.1819  # WVR!K107
$ws.Cells.Item(107, 13).Value2 = -1268.1819  # WVR!M107

# Row 126
$ws.Cells.Item(126, 8).Value2 = 4634.5386  # WVR!H126
$ws.Cells.Item(126, 10).Value2 = 7999  # WVR!J126
$ws.Cells.Item(126, 12).Value2 = 23997  # WVR!L126
$ws.Cells.Item(126, 14).Value2 = -28937  # WVR!N126
